# "Generate Report for Handback"
# Fills in the "Latest Target File" / "Latest Handback File" / "Latest Handback
# DateTime" columns on the per-locale sheets now that handback has completed,
# flips the Overview/per-locale Status from "Ready for handoff" to
# "Handed back: in sync with en-US", and widens the columns that now hold
# longer text so it stays readable.

$wb = $excel.ActiveWorkbook

$newStatus  = "Handed back: in sync with en-US"
$zhHandback = "2016-08-16 16:56:35"
$deHandback = "2016-08-16 16:56:44"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/92aa14ea99b52e60576acee20c73c298981888da/e2e/724b892e-5a63-44a9-8a22-a26a6d50ac82.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/92aa14ea99b52e60576acee20c73c298981888da/e2e/97a71edd-6aee-4803-89fa-0a1f2e37d95e.md"
$mdName1 = "724b892e-5a63-44a9-8a22-a26a6d50ac82.md"
$mdName2 = "97a71edd-6aee-4803-89fa-0a1f2e37d95e.md"

# ---------------------------------------------------------------------------
# Overview sheet: status column text + widen the two status columns (E, F)
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsOverview.Range("E1").ColumnWidth = 29.15
$wsOverview.Range("F1").ColumnWidth = 29.15

# ---------------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Status column text
$wsZh.Range("C2").Value = $newStatus
$wsZh.Range("C3").Value = $newStatus

# Column widths: Status (C), Latest Target File (I), Latest Handback File (J)
$wsZh.Range("C1").ColumnWidth = 29.15
$wsZh.Range("I1").ColumnWidth = 39.15
$wsZh.Range("J1").ColumnWidth = 39.15

# Row 2 (724b892e...)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl1, "", "", $mdName1) | Out-Null
$wsZh.Range("J2").Value = "724b892e-5a63-44a9-8a22-a26a6d50ac82.1f4387566ab9ddeb591e9375fe66ff1c5f05ac0e.zh-cn.xlf"
$wsZh.Range("K2").Value = $zhHandback

# Row 3 (97a71edd...)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl2, "", "", $mdName2) | Out-Null
$wsZh.Range("J3").Value = "97a71edd-6aee-4803-89fa-0a1f2e37d95e.8eed6d2bd3f495e26c68566cafce8bd7236a857a.zh-cn.xlf"
$wsZh.Range("K3").Value = $zhHandback

# ---------------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Status column text
$wsDe.Range("C2").Value = $newStatus
$wsDe.Range("C3").Value = $newStatus

# Column widths: Status (C), Latest Target File (I), Latest Handback File (J)
$wsDe.Range("C1").ColumnWidth = 29.15
$wsDe.Range("I1").ColumnWidth = 39.15
$wsDe.Range("J1").ColumnWidth = 39.15

# Row 2 (724b892e...)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl1, "", "", $mdName1) | Out-Null
$wsDe.Range("J2").Value = "724b892e-5a63-44a9-8a22-a26a6d50ac82.1f4387566ab9ddeb591e9375fe66ff1c5f05ac0e.de-de.xlf"
$wsDe.Range("K2").Value = $deHandback

# Row 3 (97a71edd...)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl2, "", "", $mdName2) | Out-Null
$wsDe.Range("J3").Value = "97a71edd-6aee-4803-89fa-0a1f2e37d95e.8eed6d2bd3f495e26c68566cafce8bd7236a857a.de-de.xlf"
$wsDe.Range("K3").Value = $deHandback
